$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8335
$ws1.Range("F5").Value = 35919
$ws1.Range("F8").Value = 730
$ws1.Range("F11").Value = 451
$ws1.Range("F12").Value = 819
$ws1.Range("F13").Value = 70
$ws1.Range("F14").Value = 642
$ws1.Range("F15").Value = 459
$ws1.Range("F17").Value = 588
$ws1.Range("F18").Value = 163
$ws1.Range("F19").Value = 435
$ws1.Range("F21").Value = 1129
$ws1.Range("F23").Value = 749
$ws1.Range("F24").Value = 2406
$ws1.Range("F25").Value = 895
$ws1.Range("F26").Value = 513
$ws1.Range("F27").Value = 81
$ws1.Range("F28").Value = 1106
$ws1.Range("F30").Value = 684

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 358

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 569

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 569
$ws4.Range("F3").Value = 8335
$ws4.Range("F7").Value = 35919
$ws4.Range("F10").Value = 730
$ws4.Range("F14").Value = 451
$ws4.Range("F15").Value = 358
$ws4.Range("F18").Value = 819
$ws4.Range("F19").Value = 70
$ws4.Range("F20").Value = 642
$ws4.Range("F21").Value = 459
$ws4.Range("F28").Value = 588
$ws4.Range("F29").Value = 163
$ws4.Range("F30").Value = 435
$ws4.Range("F32").Value = 1129
$ws4.Range("F34").Value = 749
$ws4.Range("F35").Value = 2406
$ws4.Range("F36").Value = 895
$ws4.Range("F37").Value = 513
$ws4.Range("F38").Value = 81
$ws4.Range("F39").Value = 1106
$ws4.Range("F42").Value = 684
